# Rename exposure variable codes to their log-transformed counterparts
# (aat1 -> ln_aat1, mpo1 -> ln_mpo1, neo1 -> ln_neo1, aat2 -> ln_aat2,
#  mpo2 -> ln_mpo2, neo2 -> ln_neo2) on the "eed_t2-dev_t23" sheet,
# column A, rows 2-19 and 32-49.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("eed_t2-dev_t23")

for ($r = 2; $r -le 19; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $current = $cell.Value2
    $cell.Value = "ln_" + $current
}

for ($r = 32; $r -le 49; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $current = $cell.Value2
    $cell.Value = "ln_" + $current
}
